$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 856, pushing existing data (old rows 856-913) down to 858-915
$ws.Rows.Item(856).Insert()
$ws.Rows.Item(856).Insert()

# Populate new row 856
$ws.Cells.Item(856, 1).Value = 5
$ws.Cells.Item(856, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(856, 3).Value = 'Maule'
$ws.Cells.Item(856, 4).Value = 44714
$ws.Cells.Item(856, 5).Value = 7
$ws.Cells.Item(856, 6).Value = 'Fruta'
$ws.Cells.Item(856, 7).Value = 100104
$ws.Cells.Item(856, 8).Value = 'Frutos de pepita'
$ws.Cells.Item(856, 9).Value = 100104002
$ws.Cells.Item(856, 10).Value = 'Manzana'
$ws.Cells.Item(856, 11).Value = 'Granny Smith'
$ws.Cells.Item(856, 12).Value = 'Primera'
$ws.Cells.Item(856, 13).Value = 200
$ws.Cells.Item(856, 14).Value = 7000
$ws.Cells.Item(856, 15).Value = 7000
$ws.Cells.Item(856, 16).Value = 7000
$ws.Cells.Item(856, 17).Value = '$/bandeja 15 kilos granel'
$ws.Cells.Item(856, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(856, 19).Value = 467
$ws.Cells.Item(856, 20).Value = 15

# Populate new row 857
$ws.Cells.Item(857, 1).Value = 5
$ws.Cells.Item(857, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(857, 3).Value = 'Maule'
$ws.Cells.Item(857, 4).Value = 44714
$ws.Cells.Item(857, 5).Value = 7
$ws.Cells.Item(857, 6).Value = 'Fruta'
$ws.Cells.Item(857, 7).Value = 100104
$ws.Cells.Item(857, 8).Value = 'Frutos de pepita'
$ws.Cells.Item(857, 9).Value = 100104002
$ws.Cells.Item(857, 10).Value = 'Manzana'
$ws.Cells.Item(857, 11).Value = 'Pink Lady'
$ws.Cells.Item(857, 12).Value = 'Primera'
$ws.Cells.Item(857, 13).Value = 250
$ws.Cells.Item(857, 14).Value = 7000
$ws.Cells.Item(857, 15).Value = 7000
$ws.Cells.Item(857, 16).Value = 7000
$ws.Cells.Item(857, 17).Value = '$/bandeja 15 kilos granel'
$ws.Cells.Item(857, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(857, 19).Value = 467
$ws.Cells.Item(857, 20).Value = 15
